# "Modificate ore in redazione relazione finale"
# Update the hours spent on the "Redazione redazione finale" (final report
# drafting) activity from 2 to 4 for all three students (rows 27-29,
# column B) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value2 = 4
$ws.Range("B28").Value2 = 4
$ws.Range("B29").Value2 = 4

$wb.Save()
